# Add the new "L6" worksheet (sheetId 7) containing per-team last-six-games summaries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "L6"

# Header row
$ws.Range("B1").Value = "Form"
$ws.Range("C1").Value = "Goals scored"
$ws.Range("D1").Value = "Goals conceded"
$ws.Range("E1").Value = "Total Goals"

# Column A holds the row numbers 1-20 as TEXT (must match shared text "1".."20", not numeric cells)
$ws.Range("A2:A21").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "1"
$ws.Cells.Item(3, 1).Value = "2"
$ws.Cells.Item(4, 1).Value = "3"
$ws.Cells.Item(5, 1).Value = "4"
$ws.Cells.Item(6, 1).Value = "5"
$ws.Cells.Item(7, 1).Value = "6"
$ws.Cells.Item(8, 1).Value = "7"
$ws.Cells.Item(9, 1).Value = "8"
$ws.Cells.Item(10, 1).Value = "9"
$ws.Cells.Item(11, 1).Value = "10"
$ws.Cells.Item(12, 1).Value = "11"
$ws.Cells.Item(13, 1).Value = "12"
$ws.Cells.Item(14, 1).Value = "13"
$ws.Cells.Item(15, 1).Value = "14"
$ws.Cells.Item(16, 1).Value = "15"
$ws.Cells.Item(17, 1).Value = "16"
$ws.Cells.Item(18, 1).Value = "17"
$ws.Cells.Item(19, 1).Value = "18"
$ws.Cells.Item(20, 1).Value = "19"
$ws.Cells.Item(21, 1).Value = "20"

# Data rows: B=Form, C=Goals scored, D=Goals conceded, E=Total Goals (team + six match values combined)
$ws.Range("B2").Value = "Atalanta,W W W D W D"
$ws.Range("C2").Value = "Atalanta,3 3 1 1 5 1"
$ws.Range("D2").Value = "Atalanta,2 2 0 1 0 1"
$ws.Range("E2").Value = "Atalanta,5 5 1 2 5 2"
$ws.Range("B3").Value = "Benevento,D L L D L L"
$ws.Range("C3").Value = "Benevento,2 0 3 2 2 0"
$ws.Range("D3").Value = "Benevento,2 1 5 2 4 2"
$ws.Range("E3").Value = "Benevento,4 1 8 4 6 2"
$ws.Range("B4").Value = "Bologna,L L W D L D"
$ws.Range("C4").Value = "Bologna,0 0 4 1 0 3"
$ws.Range("D4").Value = "Bologna,1 1 1 1 5 3"
$ws.Range("E4").Value = "Bologna,1 1 5 2 5 6"
$ws.Range("B5").Value = "Cagliari,L L W W W D"
$ws.Range("C5").Value = "Cagliari,0 0 4 1 3 1"
$ws.Range("D5").Value = "Cagliari,2 1 3 0 2 1"
$ws.Range("E5").Value = "Cagliari,2 1 7 1 5 2"
$ws.Range("B6").Value = "Crotone,L L L L W L"
$ws.Range("C6").Value = "Crotone,3 2 1 0 4 0"
$ws.Range("D6").Value = "Crotone,4 3 2 1 3 2"
$ws.Range("E6").Value = "Crotone,7 5 3 1 7 2"
$ws.Range("B7").Value = "Fiorentina,D L L W D D"
$ws.Range("C7").Value = "Fiorentina,1 2 1 2 1 3"
$ws.Range("D7").Value = "Fiorentina,1 3 3 1 1 3"
$ws.Range("E7").Value = "Fiorentina,2 5 4 3 2 6"
$ws.Range("B8").Value = "Genoa,D L L D W L"
$ws.Range("C8").Value = "Genoa,1 1 1 2 2 3"
$ws.Range("D8").Value = "Genoa,1 3 2 2 0 4"
$ws.Range("E8").Value = "Genoa,2 4 3 4 2 7"
$ws.Range("B9").Value = "Inter,W W D D W W"
$ws.Range("C9").Value = "Inter,2 1 1 1 1 2"
$ws.Range("D9").Value = "Inter,1 0 1 1 0 0"
$ws.Range("E9").Value = "Inter,3 1 2 2 1 2"
$ws.Range("B10").Value = "Juventus,W W L W D W"
$ws.Range("C10").Value = "Juventus,2 3 0 3 1 2"
$ws.Range("D10").Value = "Juventus,1 1 1 1 1 1"
$ws.Range("E10").Value = "Juventus,3 4 1 4 2 3"
$ws.Range("B11").Value = "Lazio,W W W L W W"
$ws.Range("C11").Value = "Lazio,2 1 5 2 3 4"
$ws.Range("D11").Value = "Lazio,1 0 3 5 0 3"
$ws.Range("E11").Value = "Lazio,3 1 8 7 3 7"
$ws.Range("B12").Value = "Milan,D W W L L W"
$ws.Range("C12").Value = "Milan,1 3 2 1 0 2"
$ws.Range("D12").Value = "Milan,1 1 1 2 3 0"
$ws.Range("E12").Value = "Milan,2 4 3 3 3 2"
$ws.Range("B13").Value = "Napoli,L W D W W D"
$ws.Range("C13").Value = "Napoli,1 2 1 5 2 1"
$ws.Range("D13").Value = "Napoli,2 0 1 2 0 1"
$ws.Range("E13").Value = "Napoli,3 2 2 7 2 2"
$ws.Range("B14").Value = "Parma,D L L L L L"
$ws.Range("C14").Value = "Parma,2 1 3 1 3 0"
$ws.Range("D14").Value = "Parma,2 3 4 3 4 1"
$ws.Range("E14").Value = "Parma,4 4 7 4 7 1"
$ws.Range("B15").Value = "Roma,D W L D L L"
$ws.Range("C15").Value = "Roma,2 1 1 1 2 0"
$ws.Range("D15").Value = "Roma,2 0 3 1 3 2"
$ws.Range("E15").Value = "Roma,4 1 4 2 5 2"
$ws.Range("B16").Value = "Sampdoria,D L W W L W"
$ws.Range("C16").Value = "Sampdoria,1 0 3 1 0 2"
$ws.Range("D16").Value = "Sampdoria,1 2 1 0 1 0"
$ws.Range("E16").Value = "Sampdoria,2 2 4 1 1 2"
$ws.Range("B17").Value = "Sassuolo,L W W W W D"
$ws.Range("C17").Value = "Sassuolo,1 1 3 2 1 1"
$ws.Range("D17").Value = "Sassuolo,2 0 1 1 0 1"
$ws.Range("E17").Value = "Sassuolo,3 1 4 3 1 2"
$ws.Range("B18").Value = "Spezia,L W L D L D"
$ws.Range("C18").Value = "Spezia,1 3 1 1 0 1"
$ws.Range("D18").Value = "Spezia,2 2 4 1 2 1"
$ws.Range("E18").Value = "Spezia,3 5 5 2 2 2"
$ws.Range("B19").Value = "Torino,D W W D L W"
$ws.Range("C19").Value = "Torino,2 1 3 1 0 1"
$ws.Range("D19").Value = "Torino,2 0 1 1 2 0"
$ws.Range("E19").Value = "Torino,4 1 4 2 2 1"
$ws.Range("B20").Value = "Udinese,L L W L W L"
$ws.Range("C20").Value = "Udinese,2 0 2 0 4 1"
$ws.Range("D20").Value = "Udinese,3 1 1 1 2 2"
$ws.Range("E20").Value = "Udinese,5 1 3 1 6 3"
$ws.Range("B21").Value = "Verona,W L L L L D"
$ws.Range("C21").Value = "Verona,2 0 1 1 0 1"
$ws.Range("D21").Value = "Verona,0 1 3 2 1 1"
$ws.Range("E21").Value = "Verona,2 1 4 3 1 2"

# Drop the residual NumberFormat styling so column A cells keep the default (unstyled) look, like the target
$ws.Cells.ClearFormats()

# Position the new sheet at the end, after "Goal totals v2"
$target = $wb.Worksheets.Item("Goal totals v2")
$ws.Move([ref]$null, $target)
